# Laurent coefficients computation strategy - deck update
#
# - "dérivation numérique" became "Dérivation automatique" / "dérivation automatique"
#   in the two "TextBox 1" callouts (slide 1 and slide 2).
# - The bent connector feeding each callout now points to the middle of the
#   target shape's edge (adj1 50% instead of ~2.17%).

$p = $ppt.ActivePresentation

# --- Slide 1 --------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$label1 = $s1.Shapes.Item("TextBox 1")
$label1.TextFrame.TextRange.Text = "Dérivation automatique"

$connector1 = $s1.Shapes.Item("Connector: Elbow 6")
$connector1.Adjustments.Item(1) = 0.5

# --- Slide 2 --------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$label2 = $s2.Shapes.Item("TextBox 1")
$label2.TextFrame.TextRange.Text = "dérivation automatique"

$connector2 = $s2.Shapes.Item("Connector: Elbow 4")
$connector2.Adjustments.Item(1) = 0.5
